$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("B11") "47.64"
Set-TextValue $ws.Range("C11") "3.88"
Set-TextValue $ws.Range("D11") "51.53"

Set-TextValue $ws.Range("B33") "45.47"
Set-TextValue $ws.Range("C33") "3.03"

Set-TextValue $ws.Range("C34") "39.34"
Set-TextValue $ws.Range("D34") "66.95"

Set-TextValue $ws.Range("B36") "93.58"
Set-TextValue $ws.Range("C36") "6.23"
Set-TextValue $ws.Range("D36") "99.82"

Set-TextValue $ws.Range("B40") "21.32"
Set-TextValue $ws.Range("C40") "43.38"
Set-TextValue $ws.Range("D40") "64.71"
